$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.187.82'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.832.88'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.32'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6647'
$ws.Range("E6").Value = '  -3.35%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07340'
$ws.Range("E8").Value = '  -1.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2911'
$ws.Range("E9").Value = '  -3.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.56'
$ws.Range("E10").Value = '  -3.07%  '
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").Value = '1.832.73'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.960'
$ws.Range("E13").Value = '  -2.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6636'
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.38'
$ws.Range("E15").Value = '  -4.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.064'
$ws.Range("E16").Value = '  -1.97%  '
$ws.Range("D17").Value = '29.193.62'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008249'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '224.78'
$ws.Range("E19").Value = '  -2.13%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.40'
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.116'
$ws.Range("E22").Value = '  -3.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.41'
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.604'
$ws.Range("E25").Value = '  -2.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1388'
$ws.Range("E26").Value = '  -4.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.86'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.508'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.099'
$ws.Range("E29").Value = '  -4.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.024'
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05283'
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.860'
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7473'
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.125'
$ws.Range("E35").Value = '  -0.98%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").Value = '1.310.61'
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01794'
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.714'
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9156'
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.921'
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("B42").Value = 'XinFinNetwork'
$ws.Range("C42").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.08454'
$ws.Range("E42").Value = '  +17.78%  '
$ws.Range("B43").Value = 'BabyDogeCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.00000000133'
$ws.Range("E43").Value = '  +8.92%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.88'
$ws.Range("E45").Value = '  -3.03%  '
$ws.Range("D46").Value = '1.984.14'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5164'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.760'
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.06'
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05928'
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.019'
$ws.Range("E51").Value = '  -5.39%  '
